$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Table cell paragraph: replace the "proposal" sentence with the
#    new "Desarrollado por..." sentence, splitting it into two runs
#    with a _GoBack bookmark at the split point (matches what Word
#    leaves behind after the last text edit at that spot).
# -----------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("El proyecto ha sido propuesto por la Universidad Pontificia de Salamanca a un grupo de alumnos interesados.")
$sentenceStart = $target.Start
$target.Text = "Desarrollado por alumnos de la Universidad Pontificia de Salamanca (Grado Ingenieria Informática) para la Unidad de Traumatología del Hospital Infantil Universitario Niño Jesús"

$splitPos = $sentenceStart + 83
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# -----------------------------------------------------------------
# 2) Footer: center the paragraph and drop the red tab-filler run
#    plus the trailing "Confidencial" run.
# -----------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)

$titleRng = $footer.Range
$titleRng.Find.Execute("PROYECTO FIN DE GRADO")
$titleEnd = $titleRng.End

$fullFooterEnd = $footer.Range.End

$tailRng = $footer.Range
$tailRng.SetRange($titleEnd, $fullFooterEnd - 1)
$tailRng.Text = ""

$footer.Range.Paragraphs.Item(1).Alignment = 1
